{"js": "// Office.js (Word JavaScript API)\n// Updates the worksheet date header and every three-digit-by-one-digit\n// multiplication prompt in the table, per the authoring diff.\nconst replacements = [\n  [\"2025-08-03 Sunday\", \"2025-08-04 Monday\"],\n  [\"822\u00d79=\", \"389\u00d77=\"],\n  [\"864\u00d74=\", \"266\u00d78=\"],\n  [\"829\u00d79=\", \"384\u00d79=\"],\n  [\"425\u00d77=\", \"163\u00d78=\"],\n  [\"241\u00d72=\", \"734\u00d75=\"],\n  [\"926\u00d72=\", \"154\u00d76=\"],\n  [\"242\u00d73=\", \"323\u00d79=\"],\n  [\"614\u00d72=\", \"234\u00d74=\"],\n  [\"273\u00d72=\", \"810\u00d75=\"],\n  [\"925\u00d76=\", \"401\u00d79=\"],\n  [\"717\u00d78=\", \"977\u00d79=\"],\n  [\"931\u00d76=\", \"365\u00d75=\"],\n  [\"497\u00d77=\", \"693\u00d76=\"],\n  [\"538\u00d77=\", \"636\u00d72=\"],\n  [\"862\u00d74=\", \"280\u00d73=\"],\n  [\"716\u00d79=\", \"458\u00d76=\"],\n  [\"986\u00d76=\", \"133\u00d77=\"],\n  [\"239\u00d72=\", \"867\u00d75=\"],\n  [\"579\u00d72=\", \"984\u00d79=\"],\n  [\"235\u00d79=\", \"205\u00d78=\"],\n  [\"299\u00d73=\", \"495\u00d75=\"],\n  [\"229\u00d75=\", \"939\u00d77=\"],\n  [\"886\u00d73=\", \"757\u00d72=\"],\n  [\"390\u00d74=\", \"394\u00d75=\"],\n  [\"983\u00d72=\", \"867\u00d76=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, 'Replace');\n  }\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Applies text replacements per the diff: date header + each multiplication problem cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-08-03 Sunday\", \"2025-08-04 Monday\"),\n    @(\"822\u00d79=\", \"389\u00d77=\"),\n    @(\"864\u00d74=\", \"266\u00d78=\"),\n    @(\"829\u00d79=\", \"384\u00d79=\"),\n    @(\"425\u00d77=\", \"163\u00d78=\"),\n    @(\"241\u00d72=\", \"734\u00d75=\"),\n    @(\"926\u00d72=\", \"154\u00d76=\"),\n    @(\"242\u00d73=\", \"323\u00d79=\"),\n    @(\"614\u00d72=\", \"234\u00d74=\"),\n    @(\"273\u00d72=\", \"810\u00d75=\"),\n    @(\"925\u00d76=\", \"401\u00d79=\"),\n    @(\"717\u00d78=\", \"977\u00d79=\"),\n    @(\"931\u00d76=\", \"365\u00d75=\"),\n    @(\"497\u00d77=\", \"693\u00d76=\"),\n    @(\"538\u00d77=\", \"636\u00d72=\"),\n    @(\"862\u00d74=\", \"280\u00d73=\"),\n    @(\"716\u00d79=\", \"458\u00d76=\"),\n    @(\"986\u00d76=\", \"133\u00d77=\"),\n    @(\"239\u00d72=\", \"867\u00d75=\"),\n    @(\"579\u00d72=\", \"984\u00d79=\"),\n    @(\"235\u00d79=\", \"205\u00d78=\"),\n    @(\"299\u00d73=\", \"495\u00d75=\"),\n    @(\"229\u00d75=\", \"939\u00d77=\"),\n    @(\"886\u00d73=\", \"757\u00d72=\"),\n    @(\"390\u00d74=\", \"394\u00d75=\"),\n    @(\"983\u00d72=\", \"867\u00d76=\"),\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $oldText\n    $rng.Find.Replacement.Text = $newText\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
